# Add result.txt — append a new "Perf-AIO:" benchmark section and refresh
# a couple of data points that came from a fresh fio run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh a couple of existing data points (new fio run results) -------
$ws.Range("C5").Value = 1425
$ws.Range("C12").Value = 179

# Tidy a stray trailing space in the "K IOPS " label used for the 128qd row
# of the first ("Kernel:") block.
$ws.Range("D6").Value = "K IOPS"

# --- Make room for a new "Perf-AIO:" block after the existing "Perf:" one -
# Existing rows 22-34 (the footnotes / legend block) need to shift down by
# 25 rows, to 47-59, to leave space for the new 6-row block (rows 23-28).
$ws.Rows("22:46").Insert()

# Build the new block by cloning the layout/styles of the first block
# (rows 1-6, "Kernel:") which uses the same 3-column (B/F/J) x 4-row (1qd,
# 4qd, 16qd, 128qd) shape.
$ws.Range("A1:L6").Copy($ws.Range("A23"))

# Drop the stray blank cells the clone introduces outside of the real
# layout of the new block (the template rows keep a couple of gaps).
$ws.Range("C23:L23").ClearContents()
$ws.Range("A24").ClearContents()
$ws.Range("C24:D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("I24").ClearContents()
$ws.Range("K24:L24").ClearContents()
$ws.Range("A25:A28").ClearContents()
$ws.Range("E25:E28").ClearContents()
$ws.Range("I25:I28").ClearContents()

# The new block's label column (B) isn't shaded like the cloned template's,
# so reset it back to the default style.
$ws.Range("B25:B28").Style = "Normal"

# Section title.
$ws.Range("A23").Value = "Perf-AIO:"

# New fio result figures for the Perf-AIO block.
$ws.Range("C25").Value = 109.2
$ws.Range("G25").Value = 87.2
$ws.Range("K25").Value = 73.3

$ws.Range("C26").Value = 370.1
$ws.Range("G26").Value = 321.7
$ws.Range("K26").Value = 285.8
$ws.Range("H26").Value = "K IOPS"

$ws.Range("C27").Value = 977.3
$ws.Range("G27").Value = 821.7
$ws.Range("K27").Value = 703.5

# 128qd row result is missing in the source data (only the label/units
# survive) — clear the value cells that the clone brought over.
$ws.Range("C28").ClearContents()
$ws.Range("G28").ClearContents()
$ws.Range("K28").ClearContents()

# --- Update the saved view/selection to where the author left off ---------
$ws.Range("D19").Select()
